$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1185.091
$ws.Range("I28").Value = 973.625
$ws.Range("K28").Value = 973.625
$ws.Range("M28").Value = -488.625
$ws.Range("H33").Value = 276
$ws.Range("I33").Value = 133.90475
$ws.Range("K33").Value = 133.90475
$ws.Range("M33").Value = 95.09524999999999
$ws.Range("H94").Value = 4574.875
$ws.Range("I94").Value = 4574.875
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4574.875
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -4123.875
$ws.Range("N94").ClearContents()
$ws.Range("H111").Value = 4151.2856
$ws.Range("I111").Value = 2043
$ws.Range("J111").Value = 5732.5
$ws.Range("K111").Value = 6129
$ws.Range("L111").Value = 17197.5
$ws.Range("M111").Value = -3062
$ws.Range("N111").Value = -23331.5
$ws.Range("H132").Value = 96249.28
$ws.Range("I132").Value = 103531.61
$ws.Range("K132").Value = 310594.83
$ws.Range("M132").Value = -308064.83
$ws.Range("H137").Value = 930560.25
$ws.Range("I137").Value = 2904.5293
$ws.Range("J137").Value = 1806679.5
$ws.Range("K137").Value = 8713.5879
$ws.Range("L137").Value = 5420038.5
$ws.Range("M137").Value = -6163.5879
$ws.Range("N137").Value = -5425138.5
$ws.Range("H139").Value = 107499.5
$ws.Range("J139").Value = 107499.5
$ws.Range("L139").Value = 107499.5
$ws.Range("N139").Value = -117779.5
$ws.Range("H141").Value = 2382.3333
$ws.Range("I141").Value = 2058.8
$ws.Range("K141").Value = 6176.400000000001
$ws.Range("M141").Value = -996.4000000000005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20021.14
$ws.Range("I32").Value = 20371.697
$ws.Range("J32").Value = 390
$ws.Range("K32").Value = 20371.697
$ws.Range("L32").Value = 390
$ws.Range("M32").Value = -20084.697
$ws.Range("N32").Value = -964
$ws.Range("H102").Value = 31057.834
$ws.Range("I102").Value = 37149.5
$ws.Range("J102").Value = 599.5
$ws.Range("K102").Value = 37149.5
$ws.Range("L102").Value = 599.5
$ws.Range("M102").Value = -35527.5
$ws.Range("N102").Value = -3843.5
$ws.Range("H104").Value = 56603.332
$ws.Range("J104").Value = 78905
$ws.Range("L104").Value = 78905
$ws.Range("N104").Value = -85893
$ws.Range("H110").Value = 2013.5186
$ws.Range("I110").Value = 1590.2174
$ws.Range("K110").Value = 1590.2174
$ws.Range("M110").Value = 454.7826
$ws.Range("H111").Value = 96844
$ws.Range("J111").Value = 96844
$ws.Range("L111").Value = 96844
$ws.Range("N111").Value = -105024
$ws.Range("H132").Value = 291548.75
$ws.Range("I132").Value = 318839.8
$ws.Range("J132").Value = 4992.5
$ws.Range("K132").Value = 956519.3999999999
$ws.Range("L132").Value = 14977.5
$ws.Range("M132").Value = -953989.3999999999
$ws.Range("N132").Value = -20037.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 431.85715
$ws.Range("I22").Value = 403.91666
$ws.Range("K22").Value = 403.91666
$ws.Range("M22").Value = -230.91666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 600
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H110").Value = 79975
$ws.Range("J110").Value = 79975
$ws.Range("L110").Value = 79975
$ws.Range("N110").Value = -88155
$ws.Range("H113").Value = 600
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 52905356
$ws.Range("I132").Value = 76926050
$ws.Range("J132").Value = 860520.3
$ws.Range("K132").Value = 230778150
$ws.Range("L132").Value = 2581560.9
$ws.Range("M132").Value = -230775620
$ws.Range("N132").Value = -2586620.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 2025.5
$ws.Range("J44").Value = 4001
$ws.Range("L44").Value = 12003
$ws.Range("N44").Value = -12799
$ws.Range("H55").Value = 1920.25
$ws.Range("I55").Value = 1404.9
$ws.Range("J55").Value = 4497
$ws.Range("K55").Value = 4214.700000000001
$ws.Range("L55").Value = 13491
$ws.Range("M55").Value = -4037.700000000001
$ws.Range("N55").Value = -13845
$ws.Range("H58").Value = 961.5
$ws.Range("I58").Value = 655
$ws.Range("K58").Value = 1965
$ws.Range("M58").Value = -1837
$ws.Range("H87").Value = 20013
$ws.Range("I87").Value = 14
$ws.Range("K87").Value = 42
$ws.Range("M87").Value = 1206
$ws.Range("H90").Value = 20013
$ws.Range("I90").Value = 14
$ws.Range("K90").Value = 126
$ws.Range("M90").Value = 6114
$ws.Range("H107").Value = 1231.4546
$ws.Range("J107").Value = 1996.3334
$ws.Range("L107").Value = 5989.0002
$ws.Range("N107").Value = -9829.0002
$ws.Range("H113").Value = 1250.279
$ws.Range("I113").Value = 348.75
$ws.Range("J113").Value = 1342.7435
$ws.Range("K113").Value = 1046.25
$ws.Range("L113").Value = 4028.2305
$ws.Range("M113").Value = 1123.75
$ws.Range("N113").Value = -8368.2305
$ws.Range("H115").Value = 4941
$ws.Range("I115").Value = 2414
$ws.Range("J115").Value = 9995
$ws.Range("K115").Value = 7242
$ws.Range("L115").Value = 29985
$ws.Range("M115").Value = -6067
$ws.Range("N115").Value = -32335

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 34941.293
$ws.Range("I122").Value = 53291.2
$ws.Range("J122").Value = 8727.143
$ws.Range("K122").Value = 159873.6
$ws.Range("L122").Value = 26181.429
$ws.Range("M122").Value = -157423.6
$ws.Range("N122").Value = -31081.429
$ws.Range("H141").Value = 21833.334
$ws.Range("J141").Value = 21833.334
$ws.Range("L141").Value = 21833.334
$ws.Range("N141").Value = -32193.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2774.25
$ws.Range("I61").Value = 2774.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2774.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2572.25
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 2410.8
$ws.Range("I93").Value = 2410.8
$ws.Range("K93").Value = 2410.8
$ws.Range("M93").Value = -1162.8
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H95").Value = 88275.2
$ws.Range("J95").Value = 88275.2
$ws.Range("L95").Value = 88275.2
$ws.Range("N95").Value = -93767.2
$ws.Range("H96").Value = 130000
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 89344
$ws.Range("J97").Value = 89344
$ws.Range("L97").Value = 89344
$ws.Range("N97").Value = -91326
$ws.Range("H99").Value = 33599
$ws.Range("I99").Value = 33599
$ws.Range("K99").Value = 33599
$ws.Range("M99").Value = -30604
$ws.Range("H101").Value = 119500
$ws.Range("J101").Value = 119500
$ws.Range("L101").Value = 119500
$ws.Range("N101").Value = -125990
$ws.Range("H104").Value = 80537.664
$ws.Range("J104").Value = 80537.664
$ws.Range("L104").Value = 80537.664
$ws.Range("N104").Value = -87525.664
$ws.Range("H108").Value = 100500.8
$ws.Range("J108").Value = 100500.8
$ws.Range("L108").Value = 100500.8
$ws.Range("N108").Value = -108180.8
$ws.Range("H110").Value = 82983
$ws.Range("J110").Value = 82983
$ws.Range("L110").Value = 82983
$ws.Range("N110").Value = -91163
$ws.Range("H112").Value = 60000
$ws.Range("I112").Value = 60000
$ws.Range("K112").Value = 60000
$ws.Range("M112").Value = -58523
$ws.Range("H113").Value = 2774.25
$ws.Range("I113").Value = 2774.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2774.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -604.25
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 100398
$ws.Range("J114").Value = 100398
$ws.Range("L114").Value = 100398
$ws.Range("N114").Value = -109076

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1631.6154
$ws.Range("I100").Value = 1266.375
$ws.Range("J100").Value = 2216
$ws.Range("K100").Value = 2532.75
$ws.Range("L100").Value = 4432
$ws.Range("M100").Value = -1991.75
$ws.Range("N100").Value = -5514
$ws.Range("H113").Value = 3390.9412
$ws.Range("J113").Value = 4669.1816
$ws.Range("L113").Value = 14007.5448
$ws.Range("N113").Value = -18347.5448
$ws.Range("H132").Value = 8752381
$ws.Range("I132").Value = 9584608
$ws.Range("K132").Value = 28753824
$ws.Range("M132").Value = -28751294
$ws.Range("H140").Value = 67999.5
$ws.Range("J140").Value = 67999.5
$ws.Range("L140").Value = 67999.5
$ws.Range("N140").Value = -78359.5
